$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Capture formulas/labels that need to slide left one column before any
#    writes happen (sources and destinations overlap, so read everything we
#    need up front).
# ---------------------------------------------------------------------------
$lbl_total   = $ws.Range("F9").Formula
$lbl_pctRest = $ws.Range("F11").Formula
$lbl_horasEst= $ws.Range("F13").Formula
$lbl_horas   = $ws.Range("H13").Formula
$lbl_f17     = $ws.Range("F17").Formula
$lbl_i22     = $ws.Range("I22").Formula

# ---------------------------------------------------------------------------
# 2) Plain value edit
# ---------------------------------------------------------------------------
$ws.Range("C2").Formula = 8

# ---------------------------------------------------------------------------
# 3) Copy the cell FORMATS one column to the left first (source formats are
#    still intact at this point), then overwrite the values/formulas.
# ---------------------------------------------------------------------------
$ws.Range("F9").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("H9").Copy()
$ws.Range("G9").PasteSpecial(-4122)

$ws.Range("F11").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("G11").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("H11").Copy()
$ws.Range("G11").PasteSpecial(-4122)

$ws.Range("F13").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("G13").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("H13").Copy()
$ws.Range("G13").PasteSpecial(-4122)

$ws.Range("F17").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("I22").Copy()
$ws.Range("H22").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Now write the actual values/formulas into their new homes.
# ---------------------------------------------------------------------------
$ws.Range("E9").Formula = $lbl_total
$ws.Range("F9").Formula = "=SUM(B2:B80)"
$ws.Range("G9").Formula = "=SUM(C2:C80)"

$ws.Range("E11").Formula = $lbl_pctRest
$ws.Range("F11").ClearContents()
$ws.Range("G11").Formula = "=(G9/F9)*100%"

$ws.Range("E13").Formula = $lbl_horasEst
$ws.Range("F13").Formula = "=(G9*20)/60"
$ws.Range("G13").Formula = $lbl_horas

$ws.Range("E17").Formula = $lbl_f17
$ws.Range("H22").Formula = $lbl_i22

# ---------------------------------------------------------------------------
# 5) Clear the now-vacated old cells completely (contents + formats) so they
#    disappear from the sheet the way the target workbook expects.
# ---------------------------------------------------------------------------
$ws.Range("H9").Clear()
$ws.Range("H11").Clear()
$ws.Range("H13").Clear()
$ws.Range("F17").Clear()
$ws.Range("I22").Clear()

# ---------------------------------------------------------------------------
# 6) New column D content ("Restante pro item" header / note text). Write D2
#    before D1 so the shared-string table picks up "Falta a Tabela..." before
#    "Restante pro item", matching the target workbook's string order.
# ---------------------------------------------------------------------------
$ws.Range("D2").Formula = "Falta a Tabela de Gestão do consumidor "
$ws.Range("D1").Formula = "Restante pro item"

# Give D1 a header-ish look consistent with the rest of row 1 / the summary
# box, and make sure D2 carries no special formatting.
$ws.Range("G9").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Font.Size = 18
$ws.Range("D1").Interior.Color = $ws.Range("B1").Interior.Color
$ws.Range("D2").ClearFormats()

# ---------------------------------------------------------------------------
# 7) Tidy up column widths: shift F/G/H widths to E/F/G and give D its own
#    width, mirroring the target <cols> block.
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = $ws.Columns("F").ColumnWidth
$ws.Columns("F").ColumnWidth = $ws.Columns("G").ColumnWidth
$ws.Columns("G").ColumnWidth = $ws.Columns("H").ColumnWidth
$ws.Columns("D").ColumnWidth = 37.140625

# ---------------------------------------------------------------------------
# 8) Final selection, matching the diff's new activeCell.
# ---------------------------------------------------------------------------
$ws.Range("E23").Select()
